$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.201.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.852.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.850.04'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.72%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.497.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.916.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.295.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.84%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  +5.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -3.79%  '
$ws.Range("B33").Value = 'WrappedeETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.007.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '32.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.794.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.89%  '
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("E38").Value = '  +4.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.320'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '439.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '48.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +16.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '143.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.848.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0357'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.13%  '
